$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 7-10 hold four field observations for the same visit/locality. The
# update re-sorts them: the row order becomes the reverse of what it was
# (old row10 data -> new row7, old row9 -> new row8, old row8 -> new row9,
# old row7 -> new row10). Only columns A (Id), I (Antal), Q (Ost), R (Nord),
# Z (Starttid) and AB (Sluttid) differ between the four rows, so only those
# need to be rewritten; everything else in rows 7-10 is identical already.

# New row 7  (= old row 10)
$ws.Range("A7").Value2 = 112092161
$ws.Range("I7").Value = "'10"
$ws.Range("Q7").Value2 = 584329.919621415
$ws.Range("R7").Value2 = 7048274.339291978
$ws.Range("Z7").Value = "17:22"
$ws.Range("AB7").Value = "17:22"

# New row 8  (= old row 9)
$ws.Range("A8").Value2 = 112092586
$ws.Range("I8").Value = "'20"
$ws.Range("Q8").Value2 = 584400.9675979441
$ws.Range("R8").Value2 = 7048356.949537945
$ws.Range("Z8").Value = "17:46"
$ws.Range("AB8").Value = "17:46"

# New row 9  (= old row 8) -- Antal (I) becomes blank again
$ws.Range("A9").Value2 = 112092130
$ws.Range("I9").ClearContents()
$ws.Range("Q9").Value2 = 584352.4882331375
$ws.Range("R9").Value2 = 7048231.676015709
$ws.Range("Z9").Value = "17:22"
$ws.Range("AB9").Value = "17:22"

# New row 10 (= old row 7) -- Antal (I) becomes blank again
$ws.Range("A10").Value2 = 112092066
$ws.Range("I10").ClearContents()
$ws.Range("Q10").Value2 = 584345.5636095351
$ws.Range("R10").Value2 = 7048206.515963284
$ws.Range("Z10").Value = "17:18"
$ws.Range("AB10").Value = "17:18"
